$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Imposition" / "Press Approval Task" rows for the two
# "- - 4pp Text  4p" component rows (rows 2 & 3).
$ws.Range("B2").Value = "Imposition"
$ws.Range("D2").Value = "134-Prepare files for CTP"
$ws.Range("L2").Value = "134-Prepare files for CTP"
$ws.Range("M2").Value = "134-Prepare files for CTP"

$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "169-Press Approval Task "
$ws.Range("L3").Value = "Press Approval Task"
$ws.Range("M3").Value = "169-Press Approval Task "

# Swap the "Fold" / "Webpress 1x1" rows for the "4pp Text  4p" component
# rows (rows 6 & 7).
$ws.Range("B6").Value = "Webpress 1x1"
$ws.Range("D6").Value = "370-2-D-8 5/C Sheet"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "869,015"
$ws.Range("L6").Value = "370-2-D-8 5/C Sheet"
$ws.Range("M6").Value = "311-1-D-8 5/C Sheet`n370-2-D-8 5/C Sheet"

$ws.Range("B7").Value = "Fold"
$ws.Range("D7").Value = "440-26`" Stahl 4P TD Cont"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "860,176"
$ws.Range("L7").Value = "440-26`" Stahl 4P TD Cont"
$ws.Range("M7").Value = "420-26`" Stahl 6P TD Cont`n422-26`" MBO 6P Cont`n424-26`" Stahl 4P TD Cont`n426-26`" Stahl 6P TF Cont (Glue`n428-26`" MBO 4P Cont`n432-26`" Stahl 4P TD Cont`n445-26`" Stahl USA 6P Cont. (Gl`n418-26`" Stahl 4P Pile`n440-26`" Stahl 4P TD Cont"
